# Append new trading log rows (10-15) to Sheet1, per commit:
# "Update trading results - Fri Sep 19 19:35:14 UTC 2025"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: NEAR trading attempt
$ws.Cells.Item(10, 1).Value = "2025-09-19T19:35:09.642979"
$ws.Cells.Item(10, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(10, 3).Value = "NEAR"
$ws.Cells.Item(10, 4).Value = "UNKNOWN"
$ws.Cells.Item(10, 5).Value = 3.119850135476197
$ws.Cells.Item(10, 11).Value = "ATTEMPT"
$ws.Cells.Item(10, 12).Value = "Attempting trade 1/3"

# Row 11: NEAR position opened (success)
$ws.Cells.Item(11, 1).Value = "2025-09-19T19:35:11.142328"
$ws.Cells.Item(11, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(11, 3).Value = "NEAR"
$ws.Cells.Item(11, 4).Value = "UNKNOWN"
$ws.Cells.Item(11, 5).Value = 3.119850135476197
$ws.Cells.Item(11, 6).Value = 120
$ws.Cells.Item(11, 7).Value = 1
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 11).Value = "SUCCESS"

# Row 12: SUI trading attempt
$ws.Cells.Item(12, 1).Value = "2025-09-19T19:35:11.155401"
$ws.Cells.Item(12, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(12, 3).Value = "SUI"
$ws.Cells.Item(12, 4).Value = "UNKNOWN"
$ws.Cells.Item(12, 5).Value = 3.665685532018927
$ws.Cells.Item(12, 11).Value = "ATTEMPT"
$ws.Cells.Item(12, 12).Value = "Attempting trade 2/3"

# Row 13: SUI position failed
$ws.Cells.Item(13, 1).Value = "2025-09-19T19:35:12.737999"
$ws.Cells.Item(13, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(13, 3).Value = "SUI"
$ws.Cells.Item(13, 4).Value = "UNKNOWN"
$ws.Cells.Item(13, 11).Value = "FAILED"
$ws.Cells.Item(13, 12).Value = "Trade execution failed for trade 2"

# Row 14: ADA trading attempt
$ws.Cells.Item(14, 1).Value = "2025-09-19T19:35:12.752205"
$ws.Cells.Item(14, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(14, 3).Value = "ADA"
$ws.Cells.Item(14, 4).Value = "UNKNOWN"
$ws.Cells.Item(14, 5).Value = 0.8978693919281808
$ws.Cells.Item(14, 11).Value = "ATTEMPT"
$ws.Cells.Item(14, 12).Value = "Attempting trade 3/3"

# Row 15: ADA position failed
$ws.Cells.Item(15, 1).Value = "2025-09-19T19:35:14.350933"
$ws.Cells.Item(15, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(15, 3).Value = "ADA"
$ws.Cells.Item(15, 4).Value = "UNKNOWN"
$ws.Cells.Item(15, 11).Value = "FAILED"
$ws.Cells.Item(15, 12).Value = "Trade execution failed for trade 3"
